# Daily attendance processing - 2026-01-17 11:55:49
# Reorders the "Recorded By" (column G) names for rows where the recorder
# list order changed (e.g. the real-user recorder is listed before
# "System", and mixed-case "system"/"System" duplicates are re-ordered).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of exact "before" -> "after" values for the "Recorded By" column.
$map = @{
    "System, dnasr281@gmail.com" = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System, system" = "backup@backdoor.com, system, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$colLetter = "G"

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Range("$colLetter$row")
    $current = $cell.Value2
    if ($current -ne $null -and $map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
